$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add three new rows (18, 19, 20) to the "QUÁ TRÌNH CÔNG VIỆC" table, mirroring
# the look & feel of the existing rows above them (borders + centered text).
# ---------------------------------------------------------------------------

# Pick up the formatting of existing, similarly-bordered cells so the new
# rows visually match the rest of the table (thin box border around each
# cell, horizontally + vertically centered content, and the same date
# format used elsewhere in column I).
$ws.Range("C17").Copy()
$ws.Range("C28:D28").PasteSpecial(-4122)
$ws.Range("C29:I29").PasteSpecial(-4122)
$ws.Range("C30:H30").PasteSpecial(-4122)

$ws.Range("E20").Copy()
$ws.Range("E28").PasteSpecial(-4122)

$ws.Range("F20").Copy()
$ws.Range("F28:G28").PasteSpecial(-4122)

$ws.Range("H20").Copy()
$ws.Range("H28").PasteSpecial(-4122)

$ws.Range("I14").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I30").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 28: Thảo luận chức năng 2
$ws.Range("C28").Value = 18
$ws.Range("D28").Value = "THẢO LUẬN CHỨC NĂNG 2"
$ws.Range("E28").Value = "30 PHÚT"
$ws.Range("I28").Value = 42840

# Row 29: blank separator row
$ws.Range("C29").Value = 19

# Row 30: Lập trình chức năng 2
$ws.Range("C30").Value = 20
$ws.Range("D30").Value = "LẬP TRÌNH CHỨC NĂNG 2"
$ws.Range("E30").Value = "180 phút"
$ws.Range("I30").Value = 42840

# Merge the "30 PHÚT" cell across E:H like the other duration cells above it.
$ws.Range("E28:H28").Merge()

# ---------------------------------------------------------------------------
# Update the view so the new rows are visible / selected, like in the
# uploaded workbook.
# ---------------------------------------------------------------------------
$ws.Range("D36").Select() | Out-Null
